$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = "그린 함수를 이용한 해법"
$ws.Range("E5").Value = "https://angeloyeo.github.io/2021/06/09/Greens_function.html"

$ws.Range("D6").Value = "[Python] matplotlib 으로 pandas data 그래프 그리기 :: multiple plots, Customizing Figure Layouts, scatter, boxplot"
$ws.Range("E6").Value = "https://leedakyeong.tistory.com/entry/Python-matplotlib-%EC%9C%BC%EB%A1%9C-pandas-data-%EA%B7%B8%EB%9E%98%ED%94%84-%EA%B7%B8%EB%A6%AC%EA%B8%B0-multiple-plot-Customizing-Figure-Layouts"

$ws.Range("D37").Value = "[Paper Review] MLP-Mixer: An all-MLP Architecture for Vision"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1743&mod=document&pageid=1"
